$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 2450
$ws.Range("I64").Value = 1966.6666
$ws.Range("K64").Value = 1966.6666
$ws.Range("M64").Value = -1718.6666

$ws.Range("H67").Value = 2450
$ws.Range("I67").Value = 1966.6666
$ws.Range("K67").Value = 1966.6666
$ws.Range("M67").Value = -1108.6666

$ws.Range("H76").Value = 52890.15
$ws.Range("I76").Value = 52890.15
$ws.Range("K76").Value = 52890.15
$ws.Range("M76").Value = -52575.15

$ws.Range("H79").Value = 52890.15
$ws.Range("I79").Value = 52890.15
$ws.Range("K79").Value = 52890.15
$ws.Range("M79").Value = -51798.15

$ws.Range("H116").Value = 2260
$ws.Range("I116").Value = 1600
$ws.Range("J116").Value = 2700
$ws.Range("K116").Value = 1600
$ws.Range("L116").Value = 2700
$ws.Range("M116").Value = 1842
$ws.Range("N116").Value = -9584

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H24").Value = 0
$ws.Range("J24").Value = 0
$ws.Range("L24").ClearContents()
$ws.Range("N24").Value = 0

$ws.Range("H32").Value = 14927.309
$ws.Range("I32").Value = 15408.053
$ws.Range("K32").Value = 15408.053
$ws.Range("M32").Value = -15121.053

$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("M62").ClearContents()

$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("M65").ClearContents()

$ws.Range("H100").Value = 0
$ws.Range("J100").Value = 0
$ws.Range("L100").ClearContents()
$ws.Range("N100").Value = 0

$ws.Range("H132").Value = 3532.8333
$ws.Range("I132").Value = 4175.3335
$ws.Range("J132").Value = 2462
$ws.Range("K132").Value = 12526.0005
$ws.Range("L132").Value = 7386
$ws.Range("M132").Value = -9996.000499999998
$ws.Range("N132").Value = -12446

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1637.8
$ws.Range("I86").Value = 1215.8462
$ws.Range("J86").Value = 2421.4285
$ws.Range("K86").Value = 1215.8462
$ws.Range("L86").Value = 2421.4285
$ws.Range("M86").Value = -92.84619999999995
$ws.Range("N86").Value = -4667.4285

$ws.Range("H89").Value = 1637.8
$ws.Range("I89").Value = 1215.8462
$ws.Range("J89").Value = 2421.4285
$ws.Range("K89").Value = 6079.231
$ws.Range("L89").Value = 12107.1425
$ws.Range("M89").Value = -463.2309999999998
$ws.Range("N89").Value = -23339.1425

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 823.41174
$ws.Range("I58").Value = 878.087
$ws.Range("J58").Value = 709.0909
$ws.Range("K58").Value = 878.087
$ws.Range("L58").Value = 709.0909
$ws.Range("M58").Value = -675.087
$ws.Range("N58").Value = -1115.0909

$ws.Range("H134").Value = 723.9286
$ws.Range("I134").Value = 731.5
$ws.Range("J134").Value = 678.5
$ws.Range("K134").Value = 2194.5
$ws.Range("L134").Value = 2035.5
$ws.Range("M134").Value = 340.5
$ws.Range("N134").Value = -7105.5

$ws.Range("H136").Value = 823.41174
$ws.Range("I136").Value = 878.087
$ws.Range("J136").Value = 709.0909
$ws.Range("K136").Value = 2634.261
$ws.Range("L136").Value = 2127.2727
$ws.Range("M136").Value = -84.26099999999997
$ws.Range("N136").Value = -7227.2727

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 46311336
$ws.Range("I137").Value = 25641730
$ws.Range("J137").Value = 100052310
$ws.Range("K137").Value = 76925190
$ws.Range("L137").Value = 300156930
$ws.Range("M137").Value = -76920090
$ws.Range("N137").Value = -300167130

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 28339266
$ws.Range("I70").Value = 36433356
$ws.Range("K70").Value = 36433356
$ws.Range("M70").Value = -36433086

$ws.Range("H73").Value = 28339266
$ws.Range("I73").Value = 36433356
$ws.Range("K73").Value = 36433356
$ws.Range("M73").Value = -36432420

$ws.Range("H80").Value = 4321.625
$ws.Range("I80").Value = 3605.15
$ws.Range("J80").Value = 5515.75
$ws.Range("K80").Value = 3605.15
$ws.Range("L80").Value = 5515.75
$ws.Range("M80").Value = -2607.15
$ws.Range("N80").Value = -7511.75

$ws.Range("H83").Value = 4321.625
$ws.Range("I83").Value = 3605.15
$ws.Range("J83").Value = 5515.75
$ws.Range("K83").Value = 18025.75
$ws.Range("L83").Value = 27578.75
$ws.Range("M83").Value = -13033.75
$ws.Range("N83").Value = -37562.75

$ws.Range("H122").Value = 4745
$ws.Range("I122").Value = 6879
$ws.Range("J122").Value = 2966.6667
$ws.Range("K122").Value = 20637
$ws.Range("L122").Value = 8900.000100000001
$ws.Range("M122").Value = -18187
$ws.Range("N122").Value = -13800.0001

$ws.Range("H126").Value = 2697.4
$ws.Range("I126").Value = 918
$ws.Range("J126").Value = 3587.1
$ws.Range("K126").Value = 2754
$ws.Range("L126").Value = 10761.3
$ws.Range("M126").Value = -284
$ws.Range("N126").Value = -15701.3

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 957.1429000000001
$ws.Range("I16").Value = 612.5
$ws.Range("J16").Value = 1416.6666
$ws.Range("K16").Value = 612.5
$ws.Range("L16").Value = 1416.6666
$ws.Range("M16").Value = -442.5
$ws.Range("N16").Value = -1756.6666

$ws.Range("H22").Value = 1700.75
$ws.Range("I22").Value = 1
$ws.Range("J22").Value = 2267.3333
$ws.Range("K22").Value = 1
$ws.Range("L22").Value = 2267.3333
$ws.Range("M22").Value = 294
$ws.Range("N22").Value = -2857.3333

$ws.Range("H27").Value = 1700.75
$ws.Range("I27").Value = 1
$ws.Range("J27").Value = 2267.3333
$ws.Range("K27").Value = 1
$ws.Range("L27").Value = 2267.3333
$ws.Range("M27").Value = 106
$ws.Range("N27").Value = -2481.3333

$ws.Range("H100").Value = 2600.1428
$ws.Range("I100").Value = 1500.75
$ws.Range("J100").Value = 4066
$ws.Range("K100").Value = 1500.75
$ws.Range("L100").Value = 4066
$ws.Range("M100").Value = -959.75
$ws.Range("N100").Value = -5148

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1023.0625
$ws.Range("I126").Value = 1024.6
$ws.Range("J126").Value = 1000
$ws.Range("K126").Value = 3073.8
$ws.Range("L126").Value = 3000
$ws.Range("M126").Value = -603.7999999999997
$ws.Range("N126").Value = -7940
